# "Merged cloud and info work" - reassign some Sprint1 task hours from
# "Mange" to new resource "Magnus" (halving/reducing several remaining-hour
# estimates), and move the sheet view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# --- Resource header (F10): "Mange" -> new resource "Magnus" -------------
$ws.Range("F10").Value = "Magnus"

# --- Row 19 assignment (E19 was blank) -> "Magnus", with the same cell
#     shading used for the other resource-tag cells in column E/F.
$ws.Range("F10").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Value = "Magnus"
$excel.CutCopyMode = $false

# --- Remaining-hours figures for the days that changed (columns I/J/K) ---
$ws.Range("I18:K18").Value = 1
$ws.Range("I19:K19").Value = 0.5
$ws.Range("I20:K20").Value = 0.5
$ws.Range("I23:K23").Value = 2
$ws.Range("I24:K24").Value = 1
$ws.Range("I27:K27").Value = 4

# Row 28 totals (SUM formulas) and the E34:E36 lookups on the burndown
# summary table recalculate automatically from the above.

# --- Restore view: scroll position + active selection --------------------
$ws.Activate()
$ws.Range("L27").Select()
